$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 181798.02
$ws.Cells.Item(17, 10).Value = 181798.02
$ws.Cells.Item(17, 12).Value = 545394.0599999999
$ws.Cells.Item(17, 14).Value = -545730.0599999999
$ws.Cells.Item(112, 8).Value = 1802.1428
$ws.Cells.Item(112, 10).Value = 1890
$ws.Cells.Item(112, 12).Value = 5670
$ws.Cells.Item(112, 14).Value = -7886
$ws.Cells.Item(132, 8).Value = 2281.2646
$ws.Cells.Item(132, 9).Value = 2337.516
$ws.Cells.Item(132, 10).Value = 1700
$ws.Cells.Item(132, 11).Value = 7012.548000000001
$ws.Cells.Item(132, 12).Value = 5100
$ws.Cells.Item(132, 13).Value = -4482.548000000001
$ws.Cells.Item(132, 14).Value = -10160
$ws.Cells.Item(137, 8).Value = 758.4524
$ws.Cells.Item(137, 9).Value = 598.5263
$ws.Cells.Item(137, 11).Value = 1795.5789
$ws.Cells.Item(137, 13).Value = 754.4211

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 897.5
$ws.Cells.Item(61, 9).Value = 761.1515000000001
$ws.Cells.Item(61, 11).Value = 761.1515000000001
$ws.Cells.Item(61, 13).Value = -549.1515000000001
$ws.Cells.Item(74, 8).Value = 879.4727
$ws.Cells.Item(74, 9).Value = 882.84784
$ws.Cells.Item(74, 11).Value = 882.84784
$ws.Cells.Item(74, 13).Value = -8.847840000000019
$ws.Cells.Item(77, 8).Value = 879.4727
$ws.Cells.Item(77, 9).Value = 882.84784
$ws.Cells.Item(77, 11).Value = 4414.2392
$ws.Cells.Item(77, 13).Value = -46.23919999999998
$ws.Cells.Item(132, 8).Value = 1085.1364
$ws.Cells.Item(132, 9).Value = 966.7646999999999
$ws.Cells.Item(132, 11).Value = 2900.2941
$ws.Cells.Item(132, 13).Value = -370.2941000000001
$ws.Cells.Item(136, 8).Value = 897.5
$ws.Cells.Item(136, 9).Value = 761.1515000000001
$ws.Cells.Item(136, 11).Value = 2283.4545
$ws.Cells.Item(136, 13).Value = 266.5454999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 10026.8
$ws.Cells.Item(107, 9).Value = 1175.1666
$ws.Cells.Item(107, 10).Value = 45433.332
$ws.Cells.Item(107, 11).Value = 1175.1666
$ws.Cells.Item(107, 12).Value = 45433.332
$ws.Cells.Item(107, 13).Value = 744.8334
$ws.Cells.Item(107, 14).Value = -49273.332
$ws.Cells.Item(134, 8).Value = 13965.588
$ws.Cells.Item(134, 9).Value = 1284.3188
$ws.Cells.Item(134, 11).Value = 3852.9564
$ws.Cells.Item(134, 13).Value = -1317.9564

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3516.2964
$ws.Cells.Item(31, 9).Value = 3675.7896
$ws.Cells.Item(31, 10).Value = 3137.5
$ws.Cells.Item(31, 11).Value = 3675.7896
$ws.Cells.Item(31, 12).Value = 3137.5
$ws.Cells.Item(31, 13).Value = -3380.7896
$ws.Cells.Item(31, 14).Value = -3727.5
$ws.Cells.Item(34, 8).Value = 3516.2964
$ws.Cells.Item(34, 9).Value = 3675.7896
$ws.Cells.Item(34, 10).Value = 3137.5
$ws.Cells.Item(34, 11).Value = 3675.7896
$ws.Cells.Item(34, 12).Value = 3137.5
$ws.Cells.Item(34, 13).Value = -3473.7896
$ws.Cells.Item(34, 14).Value = -3541.5
$ws.Cells.Item(58, 8).Value = 3709.1667
$ws.Cells.Item(58, 9).Value = 838.7
$ws.Cells.Item(58, 10).Value = 7297.25
$ws.Cells.Item(58, 11).Value = 838.7
$ws.Cells.Item(58, 12).Value = 7297.25
$ws.Cells.Item(58, 13).Value = -635.7
$ws.Cells.Item(58, 14).Value = -7703.25
$ws.Cells.Item(98, 8).Value = 50000
$ws.Cells.Item(98, 10).Value = 50000
$ws.Cells.Item(98, 12).Value = 50000
$ws.Cells.Item(98, 14).Value = -54492
$ws.Cells.Item(99, 8).Value = 2059.3057
$ws.Cells.Item(99, 9).Value = 1776.8182
$ws.Cells.Item(99, 10).Value = 5166.6665
$ws.Cells.Item(99, 11).Value = 1776.8182
$ws.Cells.Item(99, 12).Value = 5166.6665
$ws.Cells.Item(99, 13).Value = -278.8181999999999
$ws.Cells.Item(99, 14).Value = -8162.6665
$ws.Cells.Item(126, 8).Value = 2059.3057
$ws.Cells.Item(126, 9).Value = 1776.8182
$ws.Cells.Item(126, 10).Value = 5166.6665
$ws.Cells.Item(126, 11).Value = 5330.4546
$ws.Cells.Item(126, 12).Value = 15499.9995
$ws.Cells.Item(126, 13).Value = -2860.4546
$ws.Cells.Item(126, 14).Value = -20439.9995
$ws.Cells.Item(132, 8).Value = 1749.7164
$ws.Cells.Item(132, 9).Value = 1129.1621
$ws.Cells.Item(132, 10).Value = 2515.0667
$ws.Cells.Item(132, 11).Value = 3387.4863
$ws.Cells.Item(132, 12).Value = 7545.2001
$ws.Cells.Item(132, 13).Value = -857.4863
$ws.Cells.Item(132, 14).Value = -12605.2001
$ws.Cells.Item(134, 8).Value = 1535.4222
$ws.Cells.Item(134, 9).Value = 1632.5927
$ws.Cells.Item(134, 10).Value = 1389.6666
$ws.Cells.Item(134, 11).Value = 4897.7781
$ws.Cells.Item(134, 12).Value = 4168.9998
$ws.Cells.Item(134, 13).Value = -2362.7781
$ws.Cells.Item(134, 14).Value = -9238.9998
$ws.Cells.Item(136, 8).Value = 3709.1667
$ws.Cells.Item(136, 9).Value = 838.7
$ws.Cells.Item(136, 10).Value = 7297.25
$ws.Cells.Item(136, 11).Value = 2516.1
$ws.Cells.Item(136, 12).Value = 21891.75
$ws.Cells.Item(136, 13).Value = 33.89999999999964
$ws.Cells.Item(136, 14).Value = -26991.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(123, 8).Value = 3724.75
$ws.Cells.Item(123, 9).Value = 1225
$ws.Cells.Item(123, 11).Value = 3675
$ws.Cells.Item(123, 13).Value = -1225
$ws.Cells.Item(129, 8).Value = 43718.957
$ws.Cells.Item(129, 10).Value = 69224.336
$ws.Cells.Item(129, 12).Value = 207673.008
$ws.Cells.Item(129, 14).Value = -217673.008
$ws.Cells.Item(130, 8).Value = 2928.5715
$ws.Cells.Item(130, 9).Value = 1833.3334
$ws.Cells.Item(130, 10).Value = 3750
$ws.Cells.Item(130, 11).Value = 5500.0002
$ws.Cells.Item(130, 12).Value = 11250
$ws.Cells.Item(130, 13).Value = -480.0002000000004
$ws.Cells.Item(130, 14).Value = -21290
$ws.Cells.Item(131, 8).Value = 23954.748
$ws.Cells.Item(131, 9).Value = 144550
$ws.Cells.Item(131, 10).Value = 13402.662
$ws.Cells.Item(131, 11).Value = 433650
$ws.Cells.Item(131, 12).Value = 40207.986
$ws.Cells.Item(131, 13).Value = -428610
$ws.Cells.Item(131, 14).Value = -50287.986
$ws.Cells.Item(133, 8).Value = 4938.3335
$ws.Cells.Item(133, 10).Value = 8300
$ws.Cells.Item(133, 12).Value = 24900
$ws.Cells.Item(133, 14).Value = -35020
$ws.Cells.Item(134, 8).Value = 2772.5
$ws.Cells.Item(134, 9).Value = 1685.3846
$ws.Cells.Item(134, 10).Value = 7483.3335
$ws.Cells.Item(134, 11).Value = 5056.1538
$ws.Cells.Item(134, 12).Value = 22450.0005
$ws.Cells.Item(134, 13).Value = 13.84619999999995
$ws.Cells.Item(134, 14).Value = -32590.0005
$ws.Cells.Item(136, 8).Value = 50667.617
$ws.Cells.Item(136, 9).Value = 92670.91
$ws.Cells.Item(136, 11).Value = 278012.73
$ws.Cells.Item(136, 13).Value = -272912.73
$ws.Cells.Item(137, 8).Value = 38162.863
$ws.Cells.Item(137, 9).Value = 2126.842
$ws.Cells.Item(137, 10).Value = 106631.3
$ws.Cells.Item(137, 11).Value = 6380.526
$ws.Cells.Item(137, 12).Value = 319893.9
$ws.Cells.Item(137, 13).Value = -1280.526
$ws.Cells.Item(137, 14).Value = -330093.9
$ws.Cells.Item(138, 8).Value = 1443.3334
$ws.Cells.Item(138, 9).Value = 998.75
$ws.Cells.Item(138, 11).Value = 2996.25
$ws.Cells.Item(138, 13).Value = 2143.75
$ws.Cells.Item(139, 8).Value = 78729.46000000001
$ws.Cells.Item(139, 9).Value = 92559.09
$ws.Cells.Item(139, 10).Value = 2666.5
$ws.Cells.Item(139, 11).Value = 277677.27
$ws.Cells.Item(139, 12).Value = 7999.5
$ws.Cells.Item(139, 13).Value = -272537.27
$ws.Cells.Item(139, 14).Value = -18279.5
$ws.Cells.Item(140, 8).Value = 140193.77
$ws.Cells.Item(140, 9).Value = 276059.9
$ws.Cells.Item(140, 10).Value = 4327.636
$ws.Cells.Item(140, 11).Value = 828179.7000000001
$ws.Cells.Item(140, 12).Value = 12982.908
$ws.Cells.Item(140, 13).Value = -822999.7000000001
$ws.Cells.Item(140, 14).Value = -23342.908
$ws.Cells.Item(141, 8).Value = 8562.857
$ws.Cells.Item(141, 9).Value = 6026
$ws.Cells.Item(141, 10).Value = 9972.223
$ws.Cells.Item(141, 11).Value = 18078
$ws.Cells.Item(141, 12).Value = 29916.669
$ws.Cells.Item(141, 13).Value = -12898
$ws.Cells.Item(141, 14).Value = -40276.669

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4870.683
$ws.Cells.Item(70, 9).Value = 4864.2144
$ws.Cells.Item(70, 11).Value = 4864.2144
$ws.Cells.Item(70, 13).Value = -4594.2144
$ws.Cells.Item(73, 8).Value = 4870.683
$ws.Cells.Item(73, 9).Value = 4864.2144
$ws.Cells.Item(73, 11).Value = 4864.2144
$ws.Cells.Item(73, 13).Value = -3928.2144
$ws.Cells.Item(132, 8).Value = 1833.705
$ws.Cells.Item(132, 9).Value = 1774.0294
$ws.Cells.Item(132, 11).Value = 5322.0882
$ws.Cells.Item(132, 13).Value = -2792.0882
$ws.Cells.Item(140, 8).Value = 70980
$ws.Cells.Item(140, 10).Value = 70980
$ws.Cells.Item(140, 12).Value = 70980
$ws.Cells.Item(140, 14).Value = -81340

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1962.2909
$ws.Cells.Item(132, 9).Value = 1761.093
$ws.Cells.Item(132, 10).Value = 2683.25
$ws.Cells.Item(132, 11).Value = 5283.279
$ws.Cells.Item(132, 12).Value = 8049.75
$ws.Cells.Item(132, 13).Value = -2753.279
$ws.Cells.Item(132, 14).Value = -13109.75
$ws.Cells.Item(136, 8).Value = 1952.7843
$ws.Cells.Item(136, 9).Value = 1069.5814
$ws.Cells.Item(136, 10).Value = 6700
$ws.Cells.Item(136, 11).Value = 3208.7442
$ws.Cells.Item(136, 12).Value = 20100
$ws.Cells.Item(136, 13).Value = -658.7442000000001
$ws.Cells.Item(136, 14).Value = -25200

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 3730.125
$ws.Cells.Item(100, 9).Value = 3322
$ws.Cells.Item(100, 10).Value = 4138.25
$ws.Cells.Item(100, 11).Value = 6644
$ws.Cells.Item(100, 12).Value = 8276.5
$ws.Cells.Item(100, 13).Value = -6103
$ws.Cells.Item(100, 14).Value = -9358.5
$ws.Cells.Item(132, 8).Value = 1041.1025
$ws.Cells.Item(132, 9).Value = 971.4167
$ws.Cells.Item(132, 10).Value = 1152.6
$ws.Cells.Item(132, 11).Value = 2914.2501
$ws.Cells.Item(132, 12).Value = 3457.8
$ws.Cells.Item(132, 13).Value = -384.2501000000002
$ws.Cells.Item(132, 14).Value = -8517.799999999999
$ws.Cells.Item(136, 8).Value = 1276.56
$ws.Cells.Item(136, 9).Value = 1390.2632
$ws.Cells.Item(136, 11).Value = 4170.7896
$ws.Cells.Item(136, 13).Value = -1620.7896
